$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new expense row (row 5): item, quantity, price, source
$ws.Range("A5").Value = "Lights, Fets, cables and chargers"
$ws.Range("B5").Value = "lots!"
$ws.Range("C5").Value = 20
$ws.Range("D5").Value = "ali express"

# Match the currency number format used by the other price cells (C3/C4)
$ws.Range("C5").NumberFormat = $ws.Range("C4").NumberFormat

# Update the active selection to D6, matching the saved selection state
$ws.Range("D6").Select()
